$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '30.567.19'
$ws.Range("E2").Value = '  -0.12%  '
$ws.Range("D3").Value = '1.922.63'
$ws.Range("E3").Value = '  -0.03%  '
$ws.Range("E4").Value = '  -0.01%  '
$ws.Range("D5").Value = "'245.59"
$ws.Range("E5").Value = '  -1.16%  '
$ws.Range("E6").Value = '  -0.03%  '
$ws.Range("D7").Value = "'0.4832"
$ws.Range("E7").Value = '  +2.05%  '
$ws.Range("D8").Value = "'0.2893"
$ws.Range("E8").Value = '  -0.48%  '
$ws.Range("D9").Value = "'0.06794"
$ws.Range("E9").Value = '  -1.17%  '
$ws.Range("D10").Value = "'112.29"
$ws.Range("E10").Value = '  +6.54%  '
$ws.Range("D11").Value = "'19.47"
$ws.Range("E11").Value = '  +5.86%  '
$ws.Range("D12").Value = '1.918.53'
$ws.Range("E12").Value = '  -0.25%  '
$ws.Range("D13").Value = "'0.07572"
$ws.Range("D14").Value = "'5.460"
$ws.Range("E14").Value = '  +2.09%  '
$ws.Range("D15").Value = "'0.6744"
$ws.Range("E15").Value = '  +0.61%  '
$ws.Range("E16").Value = '  +1.70%  '
$ws.Range("D17").Value = '30.540.69'
$ws.Range("E17").Value = '  -0.24%  '
$ws.Range("D18").Value = "'0.000007671"
$ws.Range("E18").Value = '  +0.35%  '
$ws.Range("D19").Value = "'13.02"
$ws.Range("E19").Value = '  +0.64%  '
$ws.Range("D20").Value = "'1.0000"
$ws.Range("D21").Value = "'5.501"
$ws.Range("E21").Value = '  -0.82%  '
$ws.Range("D22").Value = '2.163.66'
$ws.Range("E22").Value = '  -0.44%  '
$ws.Range("D23").Value = "'0.9998"
$ws.Range("E23").Value = '  -0.08%  '
$ws.Range("D24").Value = "'6.449"
$ws.Range("E24").Value = '  +0.15%  '
$ws.Range("D25").Value = "'9.473"
$ws.Range("E25").Value = '  -0.35%  '
$ws.Range("E26").Value = '  -0.48%  '
$ws.Range("D27").Value = "'20.28"
$ws.Range("E27").Value = '  -2.28%  '
$ws.Range("D28").Value = "'2.100"
$ws.Range("E28").Value = '  -1.22%  '
$ws.Range("D29").Value = "'0.1067"
$ws.Range("E29").Value = '  -0.46%  '
$ws.Range("E30").Value = '  +2.29%  '
$ws.Range("D31").Value = "'4.137"
$ws.Range("E31").Value = '  -1.11%  '
$ws.Range("D32").Value = "'4.056"
$ws.Range("E32").Value = '  +0.00%  '
$ws.Range("D33").Value = "'0.04941"
$ws.Range("E33").Value = '  -1.61%  '
$ws.Range("D34").Value = "'0.7358"
$ws.Range("E34").Value = '  +0.26%  '
$ws.Range("D35").Value = "'1.138"
$ws.Range("E35").Value = '  -0.79%  '
$ws.Range("E36").Value = '  -0.61%  '
$ws.Range("D37").Value = "'0.02026"
$ws.Range("E37").Value = '  -2.37%  '
$ws.Range("D38").Value = "'2.691"
$ws.Range("E38").Value = '  +0.14%  '
$ws.Range("D39").Value = "'2.019"
$ws.Range("E39").Value = '  -1.25%  '
$ws.Range("D40").Value = "'109.62"
$ws.Range("E40").Value = '  -1.49%  '
$ws.Range("D41").Value = "'0.4430"
$ws.Range("E41").Value = '  -0.13%  '
$ws.Range("D42").Value = "'0.8686"
$ws.Range("E42").Value = '  -0.92%  '
$ws.Range("D43").Value = "'5.824"
$ws.Range("E43").Value = '  -1.21%  '
$ws.Range("E44").Value = '  +0.02%  '
$ws.Range("D45").Value = "'69.22"
$ws.Range("E45").Value = '  +2.19%  '
$ws.Range("D46").Value = "'7.247"
$ws.Range("E46").Value = '  -0.83%  '
$ws.Range("D47").Value = "'48.78"
$ws.Range("E47").Value = '  +1.96%  '
$ws.Range("D48").Value = "'9.204"
$ws.Range("E48").Value = '  -2.10%  '
$ws.Range("D49").Value = "'0.1231"
$ws.Range("E49").Value = '  -1.20%  '
$ws.Range("D50").Value = "'34.81"
$ws.Range("E50").Value = '  -0.67%  '
$ws.Range("E51").Value = '  -0.26%  '
